# camparam.xlsx — "added hdf5 and fits video reader"
#
# The canonical-XML diff for this commit touches camparam.xlsx only in
# ways that are artifacts of the workbook having been re-saved by
# LibreOffice Calc alongside the (unrelated, code-only) hdf5/fits reader
# change:
#   - workbookView tabRatio (984 -> 989) and a new <extLst><loext:extCalcPr.../>
#     block: LibreOffice-specific view/recalc bookkeeping written on every
#     save, not reachable through the Excel object model.
#   - the two <col> width values shrink very slightly: LibreOffice's
#     recomputed default column width for the save environment's font
#     metrics, not a user-set width.
#   - a new shared string "farneback" appears in xl/sharedStrings.xml,
#     shifting every later <si> index by +1 — but every cell that used to
#     reference "hs" (shared-string index 2) is re-pointed to shared-string
#     index 3, which is still "hs" after the insertion. No cell's value,
#     formula, style, comment, or anything else user-visible actually
#     changes (confirmed by diffing before/after). "farneback" was already
#     listed as a valid ofmethod option in A2's comment; it lands in the
#     shared-string table unused/unreferenced by any cell, which only
#     happens as a side effect of LibreOffice's save process (e.g. an
#     autocomplete/undo cache), not something settable via
#     Range.Value/Formula on the Excel COM object model.
#
# So there is no cell content, formula, formatting, or structural edit to
# make here: the workbook's data is byte-for-byte equivalent before and
# after at the object-model level. We simply touch/re-save the workbook
# through the COM surface without modifying any range, matching the
# no-op-on-content nature of the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-select the same cell the sheet was already on; no values, formulas,
# or formatting are modified anywhere in the workbook.
$ws.Range("B2").Select()

$wb.Save()
